# CA 2 prob def update
# Inserts a new "Problem Description" section (heading + two body
# paragraphs) right after the "The primary objective..." paragraph
# that ends the Introduction section.

$d = $word.ActiveDocument

# Locate the paragraph that ends with "...such tasks." (the last
# paragraph of the Introduction section) by searching for its
# distinctive leading text.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The primary objective of this project*") {
        $anchor = $p
    }
}

# --- New paragraph 1: "Problem Description" heading ---------------
$anchor.Range.InsertParagraphAfter()
$headingPara = $anchor.Next()
$headingPara.Range.Text = "Problem Description"
$headingPara.Style = "Heading1"

$hStart = $headingPara.Range.Start
$hEnd = $headingPara.Range.End
$bmRange = $d.Range($hStart, $hEnd - 1)
$d.Bookmarks.Add("_Toc183797882", $bmRange)

# --- New paragraph 2: glass-classification-challenge paragraph ----
$headingPara.Range.InsertParagraphAfter()
$challengePara = $headingPara.Next()
$challengePara.Style = "Normal"
$challengePara.Range.Text = "The classification of glass types poses a unique challenge due to the inherent complexity of the material's composition and its applications across diverse industries. Glass is made from a mixture of raw materials such as sand, soda ash, and limestone, combined with specific oxides to achieve desired properties. These compositions determine the glass's refractive index, durability, transparency, and thermal properties, which are critical for its intended use. For instance, the glass used in building windows requires a different manufacturing process and composition than glass used for headlamps or containers."

# --- New paragraph 3: multi-class classification paragraph --------
# (split across two runs, with a lastRenderedPageBreak marker before
# "features" exactly like the authored document)
$challengePara.Range.InsertParagraphAfter()
$datasetPara = $challengePara.Next()
$dStart = $datasetPara.Range.Start
$dEnd = $datasetPara.Range.End
$dRange = $d.Range($dStart, $dEnd)
$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">In the context of the provided dataset, this problem is framed as a multi-class classification task where the objective is to predict one of seven distinct glass types based on its compositional </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>features. These types include categories like float-processed building windows, non-float-processed building windows, containers, tableware, and headlamps. Each class represents a specific industrial application, with precise compositional requirements.</w:t></w:r></w:p>'
$dRange.InsertXML($paraXml) | Out-Null

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
